# LCRY-1700.docx template update
#
# 1) The document's single table was pinned to a fixed width (9907 dxa).
#    Switch it to automatic sizing (w:tblW w:w="0" w:type="auto").
# 2) The page size was A4 (11906 x 16838 dxa). Switch the section to
#    US Letter (12240 x 15840 dxa).

$d = $word.ActiveDocument

# --- 1. Table width: fixed -> automatic -------------------------------
$table = $d.Tables.Item(1)
$table.PreferredWidthType = 1   # wdPreferredWidthAuto
$table.PreferredWidth = 0

# --- 2. Page size: A4 -> US Letter -------------------------------------
# 12240 dxa / 20 = 612 pt (width), 15840 dxa / 20 = 792 pt (height)
$d.PageSetup.PageWidth = 612
$d.PageSetup.PageHeight = 792
